$d = $word.ActiveDocument

function Find-ParagraphContaining($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $para = $d.Paragraphs.Item($i)
        if ($para.Range.Text.Contains($needle)) {
            return $para
        }
    }
    return $null
}

# --- 1. Fix the VML image height (352.5pt x 177.85pt -> 352.5pt x 177.75pt) ---
# This legacy VML w:pict picture is not exposed via InlineShapes/Shapes in the
# object model (InlineShapes.Count is 0 for this document), so replace its
# containing paragraph's XML directly -- same content, only the height
# attribute corrected -- via Range.InsertXML.
$picturePara = $d.Paragraphs.Item(12)
$pictureFragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" w:rsidR="008E6AD2" w:rsidRDefault="00080D5C" w:rsidP="009634B4"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Comic Sans MS" w:eastAsia="Times New Roman" w:hAnsi="Comic Sans MS" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="en-CA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Comic Sans MS" w:eastAsia="Times New Roman" w:hAnsi="Comic Sans MS" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="en-CA"/></w:rPr><w:pict><v:shapetype id="_x0000_t75" coordsize="21600,21600" o:spt="75" o:preferrelative="t" path="m@4@5l@4@11@9@11@9@5xe" filled="f" stroked="f"><v:stroke joinstyle="miter"/><v:formulas><v:f eqn="if lineDrawn pixelLineWidth 0"/><v:f eqn="sum @0 1 0"/><v:f eqn="sum 0 0 @1"/><v:f eqn="prod @2 1 2"/><v:f eqn="prod @3 21600 pixelWidth"/><v:f eqn="prod @3 21600 pixelHeight"/><v:f eqn="sum @0 0 1"/><v:f eqn="prod @6 1 2"/><v:f eqn="prod @7 21600 pixelWidth"/><v:f eqn="sum @8 21600 0"/><v:f eqn="prod @7 21600 pixelHeight"/><v:f eqn="sum @10 21600 0"/></v:formulas><v:path o:extrusionok="f" gradientshapeok="t" o:connecttype="rect"/><o:lock v:ext="edit" aspectratio="t"/></v:shapetype><v:shape id="_x0000_i1025" type="#_x0000_t75" style="width:352.5pt;height:177.75pt"><v:imagedata r:id="rId7" o:title="Avicii"/></v:shape></w:pict></w:r></w:p>'
$picturePara.Range.InsertXML($pictureFragment) | Out-Null

# --- 2. Remove the _GoBack bookmark that sits after "Visual Composition: " ---
$visualCompositionPara = Find-ParagraphContaining("Visual Composition:")
$visualFragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00080D5C" w:rsidRPr="00411346" w:rsidRDefault="00080D5C" w:rsidP="00411346"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Comic Sans MS" w:eastAsia="Times New Roman" w:hAnsi="Comic Sans MS" w:cs="Arial"/><w:b/><w:color w:val="000000"/><w:lang w:eastAsia="en-CA"/></w:rPr></w:pPr><w:r w:rsidRPr="00411346"><w:rPr><w:rFonts w:ascii="Comic Sans MS" w:eastAsia="Times New Roman" w:hAnsi="Comic Sans MS" w:cs="Arial"/><w:b/><w:color w:val="000000"/><w:lang w:eastAsia="en-CA"/></w:rPr><w:t xml:space="preserve">Visual Composition: </w:t></w:r></w:p>'
$visualCompositionPara.Range.InsertXML($visualFragment) | Out-Null

# --- 3. Split "Choosing colours that are analogous and not too similar" into
#        "Choose" + " colours that are analogous and not too similar" (the
#        second run switches to Comic Sans MS) and move the _GoBack bookmark
#        to the end of this paragraph ---
$coloursPara = Find-ParagraphContaining("Choosing colours that are analogous")
$coloursFragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00080D5C" w:rsidRPr="00411346" w:rsidRDefault="00080D5C" w:rsidP="00411346"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Comic Sans MS" w:eastAsia="Times New Roman" w:hAnsi="Comic Sans MS" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-CA"/></w:rPr></w:pPr><w:r w:rsidRPr="00411346"><w:rPr><w:rFonts w:ascii="Comic Sans MS" w:eastAsia="Times New Roman" w:hAnsi="Comic Sans MS" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-CA"/></w:rPr><w:t>Choose</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Comic Sans MS" w:eastAsia="Times New Roman" w:hAnsi="Comic Sans MS" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-CA"/></w:rPr><w:t xml:space="preserve"> colours that are analogous and not too similar</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$coloursPara.Range.InsertXML($coloursFragment) | Out-Null

Write-Host "Edits applied."
